$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits around "dopo" in the git-commit
#    paragraph. The edit relocates it to a brand-new empty paragraph near the
#    end of the document, so remove it from its old spot first.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# ---------------------------------------------------------------------------
# 2) Add the new "GitHub branch" paragraph plus the relocated bookmark
#    paragraph. They both land where the old lone empty paragraph
#    (right after "...che voglio segnalare.") used to be.
# ---------------------------------------------------------------------------
$targetIndex = 38
$emptyPara = $d.Paragraphs($targetIndex).Range.Duplicate

# Split the single empty paragraph into two empty paragraphs: the first will
# hold the new text, the second will hold the relocated bookmark.
$emptyPara.InsertParagraphAfter()

# --- fill the first (now-target) paragraph with the new text ---
$target = $d.Paragraphs($targetIndex).Range.Duplicate
$startPos = $target.Start

$introText = "In GitHub quando creo un nuovo Repository verrà automaticamente assegnato il branch principale chiamato master. Se voglio che il mio branch di default sia diverso dal master, basta andare sul setting del mio Repository, sulla sezione branch, rinominare il mio branch di default o cambiarlo (se sono già stati creati altri branch). "
$boldText = "Ricorda che se vai modificare questa opzione, potresti perdere i file o altro."

$target.InsertBefore($introText + $boldText)

$boldStart = $startPos + $introText.Length
$boldEnd = $boldStart + $boldText.Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = 1

# --- add the relocated bookmark into the second (still empty) paragraph ---
$bmParaIndex = $targetIndex + 1
$bmPara = $d.Paragraphs($bmParaIndex).Range.Duplicate

# A bookmark collapsed exactly on an empty paragraph's boundary gets dropped,
# so temporarily insert a couple of placeholder characters, bookmark the
# point between them, then delete the placeholders from both sides of the
# bookmark (never crossing it in a single delete).
$bmPara.InsertBefore("ZZ")
$bmParaRange = $d.Paragraphs($bmParaIndex).Range
$midPos = $bmParaRange.Start + 1

$bmRange = $d.Range($midPos, $midPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$beforeChar = $d.Range($bmParaRange.Start, $midPos)
$beforeChar.Delete()

$bmParaRange2 = $d.Paragraphs($bmParaIndex).Range
$afterChar = $d.Range($bmParaRange2.Start, $bmParaRange2.End - 1)
$afterChar.Delete()
